$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks numeric need to stay as text,
# matching the source data which stores them as inline strings.
# Force text format before assignment so Excel does not auto-convert them to numbers.

$ws.Range("D2").Value = '30.276.72'
$ws.Range("E2").Value = '  +5.15%  '

$ws.Range("D3").Value = '1.913.85'
$ws.Range("E3").Value = '  +5.56%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.14'
$ws.Range("E5").Value = '  +1.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5157'
$ws.Range("E7").Value = '  +3.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '45.77'
$ws.Range("E8").Value = '  +6.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2987'
$ws.Range("E9").Value = '  +7.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06799'
$ws.Range("E10").Value = '  +6.59%  '

$ws.Range("D11").Value = '1.914.40'
$ws.Range("E11").Value = '  +5.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '17.44'
$ws.Range("E12").Value = '  +4.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07380'
$ws.Range("E13").Value = '  +3.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6984'
$ws.Range("E14").Value = '  +7.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.68'
$ws.Range("E15").Value = '  +7.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.886'
$ws.Range("E16").Value = '  +3.78%  '

$ws.Range("D17").Value = '30.287.82'
$ws.Range("E17").Value = '  +5.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008048'
$ws.Range("E18").Value = '  +9.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9996'
$ws.Range("E19").Value = '  +0.11%  '

$ws.Range("E20").Value = '  +6.09%  '

$ws.Range("D21").Value = '2.161.21'
$ws.Range("E21").Value = '  +5.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9988'
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.849'
$ws.Range("E23").Value = '  +5.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.746'
$ws.Range("E24").Value = '  +7.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.172'
$ws.Range("E25").Value = '  +3.17%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.51'
$ws.Range("E26").Value = '  +1.91%  '

$ws.Range("B27").Value = 'BitcoinCash'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '138.65'
$ws.Range("E27").Value = '  +19.76%  '

$ws.Range("E28").Value = '  +7.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.017'
$ws.Range("E29").Value = '  +6.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.399'
$ws.Range("E30").Value = '  +0.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.264'
$ws.Range("E31").Value = '  +2.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08833'
$ws.Range("E32").Value = '  +5.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.016'
$ws.Range("E33").Value = '  +4.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05112'
$ws.Range("E34").Value = '  +2.79%  '

$ws.Range("E35").Value = '  +6.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7218'
$ws.Range("E36").Value = '  +6.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687'
$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.842'
$ws.Range("E38").Value = '  +3.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.315'
$ws.Range("E39").Value = '  +5.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9767'
$ws.Range("E40").Value = '  +0.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01697'
$ws.Range("E41").Value = '  +6.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.092'
$ws.Range("E42").Value = '  +1.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '105.97'
$ws.Range("E43").Value = '  +4.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4323'
$ws.Range("E44").Value = '  +5.04%  '

$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("E46").Value = '  +6.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1283'
$ws.Range("E47").Value = '  +4.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05761'
$ws.Range("E48").Value = '  +4.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.33'
$ws.Range("E49").Value = '  +5.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.479'
$ws.Range("E50").Value = '  +3.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3831'
$ws.Range("E51").Value = '  +4.91%  '
